# Apply updated crypto price/volume figures (and the MultiversX/Algorand row swap)
# to match the refreshed GitHub Actions data pull.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'43.906.47"
$ws.Range("E2").Value = "'  +0.46%  "
$ws.Range("D3").Value = "'2.297.71"
$ws.Range("E3").Value = "'  +0.17%  "
$ws.Range("E4").Value = "'  +0.08%  "
$ws.Range("D5").Value = "'112.71"
$ws.Range("E5").Value = "'  +16.80%  "
$ws.Range("D6").Value = "'269.90"
$ws.Range("E6").Value = "'  +0.20%  "
$ws.Range("D7").Value = "'0.627"
$ws.Range("E7").Value = "'  +0.61%  "
$ws.Range("E8").Value = "'  +0.28%  "
$ws.Range("D9").Value = "'0.622"
$ws.Range("E9").Value = "'  +2.32%  "
$ws.Range("D10").Value = "'48.07"
$ws.Range("E10").Value = "'  +6.10%  "
$ws.Range("D11").Value = "'0.0952"
$ws.Range("E11").Value = "'  +1.90%  "
$ws.Range("D12").Value = "'9.17"
$ws.Range("E12").Value = "'  +16.42%  "
$ws.Range("E13").Value = "'  +0.38%  "
$ws.Range("D14").Value = "'15.91"
$ws.Range("E14").Value = "'  +1.35%  "
$ws.Range("D15").Value = "'2.640.89"
$ws.Range("E15").Value = "'  +0.10%  "
$ws.Range("D16").Value = "'0.855"
$ws.Range("E16").Value = "'  +0.23%  "
$ws.Range("D17").Value = "'2.287.12"
$ws.Range("E17").Value = "'  -0.59%  "
$ws.Range("D18").Value = "'43.760.26"
$ws.Range("E18").Value = "'  +0.11%  "
$ws.Range("D19").Value = "'0.0000111"
$ws.Range("E19").Value = "'  -1.72%  "
$ws.Range("D20").Value = "'6.78"
$ws.Range("E20").Value = "'  +9.45%  "
$ws.Range("D21").Value = "'72.38"
$ws.Range("E21").Value = "'  +0.38%  "
$ws.Range("E22").Value = "'  -3.28%  "
$ws.Range("D23").Value = "'232.90"
$ws.Range("E23").Value = "'  +0.03%  "
$ws.Range("D24").Value = "'9.80"
$ws.Range("E24").Value = "'  +7.62%  "
$ws.Range("D25").Value = "'2.82"
$ws.Range("E25").Value = "'  +4.90%  "
$ws.Range("E26").Value = "'  +0.03%  "
$ws.Range("E27").Value = "'  +3.98%  "
$ws.Range("D28").Value = "'42.22"
$ws.Range("E28").Value = "'  +9.70%  "
$ws.Range("E29").Value = "'  -2.07%  "
$ws.Range("E30").Value = "'  -0.51%  "
$ws.Range("D31").Value = "'175.89"
$ws.Range("E31").Value = "'  +0.62%  "
$ws.Range("D32").Value = "'0.0940"
$ws.Range("E32").Value = "'  +4.22%  "
$ws.Range("D33").Value = "'21.59"
$ws.Range("E33").Value = "'  -1.18%  "
$ws.Range("E34").Value = "'  +4.91%  "
$ws.Range("E35").Value = "'  +1.57%  "
$ws.Range("D36").Value = "'4.70"
$ws.Range("E36").Value = "'  +3.99%  "
$ws.Range("D37").Value = "'0.0365"
$ws.Range("E37").Value = "'  +3.74%  "
$ws.Range("D38").Value = "'0.108"
$ws.Range("E38").Value = "'  +1.34%  "
$ws.Range("D39").Value = "'3.81"
$ws.Range("E39").Value = "'  +11.84%  "
$ws.Range("B40").Value = "Algorand"
$ws.Range("C40").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D40").Value = "'0.244"
$ws.Range("E40").Value = "'  +2.47%  "
$ws.Range("B41").Value = "MultiversX"
$ws.Range("C41").Value = "https://coinranking.com/coin/omwkOTglq+multiversx-egld"
$ws.Range("D41").Value = "'73.95"
$ws.Range("E41").Value = "'  +14.68%  "
$ws.Range("D42").Value = "'2.40"
$ws.Range("E42").Value = "'  +3.37%  "
$ws.Range("D43").Value = "'13.68"
$ws.Range("E43").Value = "'  +12.59%  "
$ws.Range("D44").Value = "'6.41"
$ws.Range("E44").Value = "'  +23.89%  "
$ws.Range("E45").Value = "'  +0.08%  "
$ws.Range("E46").Value = "'  +4.39%  "
$ws.Range("D47").Value = "'8.80"
$ws.Range("E47").Value = "'  +0.42%  "
$ws.Range("D48").Value = "'102.70"
$ws.Range("E48").Value = "'  +5.47%  "
$ws.Range("D49").Value = "'0.0999"
$ws.Range("E49").Value = "'  -1.84%  "
$ws.Range("E50").Value = "'  +3.63%  "
$ws.Range("D51").Value = "'0.468"
$ws.Range("E51").Value = "'  +7.74%  "
